$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list values (Price and Volume(1h) columns) per upstream data refresh.
# Force target cells to remain Text-typed (matches source workbook storing these as
# inline/shared strings, e.g. "306.80" not the number 306.8) by applying a Text number
# format before writing the literal string value.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "306.80"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.17%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.27"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.41%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.099"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.98%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07606"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.21%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.607"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.77%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "2.488"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.45%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9013"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.05%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1118"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "12.74%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1764"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.48%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09204"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.60%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04311"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.51%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.61%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001253"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.93%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005799"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.36%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.360"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.30%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.253"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.01%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.68%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.576"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-6.15%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1351"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.11%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2682"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-10.59%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04201"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.66%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.13%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004070"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.14%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.63%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.94%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02388"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "1.19%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05178"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.22%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007777"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.06%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.20%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006958"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.66%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.12%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008536"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "15.18%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3046"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-8.23%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006442"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.62%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.07%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-11.96%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.01028"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "228.15%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.07%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.07%"
